# Fix spelling of "TRENTINO-ALTO ADIGE/S_DTIROL" -> "TRENTINO-ALTO ADIGE/SÜDTIROL"
# (Art. 84 commi 1 e 2)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A59").Value = "TRENTINO-ALTO ADIGE/SÜDTIROL"
$ws.Range("B59").Value = "TRENTINO-ALTO ADIGE/SÜDTIROL - 01"

# Match the author's saved view state: scrolled so row 49 is at the top,
# with cell B59 selected.
$ws.Application.ActiveWindow.ScrollRow = 49
$ws.Range("B59").Select()
